$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new study-set row (row 6) below the existing data table.
$ws.Range("A6").Value = "2023-09-28 21:41:51 7_5_8829598"
$ws.Range("B6").Value = "dddfdss"
$ws.Range("C6").Value = "sdfsd"
